# Append a new sentence about the postorder traversal right after the
# existing "inorder traversal..." sentence, before the trailing
# _GoBack bookmark, matching the target run/proofErr structure exactly.

$d = $word.ActiveDocument

# Locate the end of the sentence we are appending after.
$anchor = $d.Content
$found = $anchor.Find.Execute("inorder traversal. To do so, start by visiting each node of the left subtree, then the root node, then each node of the right subtree, and append each node to a list L.")
if (-not $found) {
    throw "Could not find the inorder-traversal sentence to anchor the insertion."
}
$anchor.Collapse(0)

# Insert a one-character placeholder right at the insertion point; it
# gives InsertXML a concrete (non-empty) range to replace in place, and
# it sits inside the paragraph rather than at the absolute end of the
# document, which keeps the new content ahead of the _GoBack bookmark.
$anchor.InsertBefore("X")
$target = $d.Range($anchor.Start, $anchor.Start + 1)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> I then implemented the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/></w:rPr><w:t>postorder</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> traversal by visiting the left subtree, then the right subtree, then the root node, appending each visited node to a list L as the tree is traversed.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.InsertXML($xml)
